$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New tracker data rows for 2025-09-12 (serial date 45912)
$goals = @(
    @{ Id = "G2"; Name = "Workout" },
    @{ Id = "G3"; Name = "Eat Healthy" },
    @{ Id = "G4"; Name = "Read Book" },
    @{ Id = "G5"; Name = "Investment Plan" },
    @{ Id = "G6"; Name = "Spend 10 Hours without phone" }
)

$startRow = 27
$dateSerial = 45912
$progress = 0.9514656876067488
$percentage = 0
$change = -0.01

for ($i = 0; $i -lt $goals.Count; $i++) {
    $row = $startRow + $i
    $goal = $goals[$i]

    $ws.Cells.Item($row, 1).Value = $goal.Id
    $ws.Cells.Item($row, 2).Value = $goal.Name
    $ws.Cells.Item($row, 3).Value = $dateSerial
    $ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($row, 4).Value = $progress
    $ws.Cells.Item($row, 5).Value = $percentage
    $ws.Cells.Item($row, 6).Value = $change
}
